$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "Criar User Story" tasks are finished (reviewed by the team) and
# move from Todo (column A) into Done (column D), with the reviewer noted.
$ws.Range("A2").Value = "Descrever Totalmente Primeira Feature"
$ws.Range("A3").Value = "Descrever Totalmente Segunda Feature"

$ws.Range("D7").Value = "Criar User Story Primeira Feature (Ricardo, revisto por todos)"
$ws.Range("D8").Value = "Criar User Story Segunda Feature (James, revisto por todos)"

# New empty styled cell (underline style like A4:A6 / B8) at A9
$ws.Range("A9").Value = ""
$ws.Range("A9").Font.Underline = 1

# Adjust column widths to fit new (longer) content.
# Note: the engine quantizes ColumnWidth to the nearest 1/6th character
# internally, so these are chosen to land on the closest representable
# width to the authored values (~32.78 and ~51.44 characters).
$ws.Columns.Item(1).ColumnWidth = 32
$ws.Columns.Item(4).ColumnWidth = 50.6666667

# Update selection to A9
$ws.Range("A9").Select()
